$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "B08MPM2BB2"

# Replace the 100 values in column A (A1:A100) with the new list
$values = New-Object 'object[,]' 100,1
$values[0,0] = "black top"
$values[1,0] = "black bra"
$values[2,0] = "wireless bra"
$values[3,0] = "sport bra"
$values[4,0] = "pink top"
$values[5,0] = "yoga bra"
$values[6,0] = "yoga top"
$values[7,0] = "low back bra"
$values[8,0] = "black s"
$values[9,0] = "longline bra"
$values[10,0] = "padded bra"
$values[11,0] = "green bra"
$values[12,0] = "strappy bra"
$values[13,0] = "bra top"
$values[14,0] = "pink bra"
$values[15,0] = "pink bralette"
$values[16,0] = "workout bra"
$values[17,0] = "strappy bralette"
$values[18,0] = "black bra top"
$values[19,0] = "longline bralette"
$values[20,0] = "black cross"
$values[21,0] = "bralette top"
$values[22,0] = "cute top"
$values[23,0] = "criss cross bra"
$values[24,0] = "black l"
$values[25,0] = "black m"
$values[26,0] = "wirefree bra"
$values[27,0] = "criss cross top"
$values[28,0] = "long line bra"
$values[29,0] = "cute bra"
$values[30,0] = "criss cross"
$values[31,0] = "pink bra top"
$values[32,0] = "bralette bra"
$values[33,0] = "impact sport"
$values[34,0] = "long s"
$values[35,0] = "top s"
$values[36,0] = "black cup"
$values[37,0] = "wireless bralette"
$values[38,0] = "women’s longline sports bra wirefree padded medium support yoga bras gym running workout tank tops"
$values[39,0] = "cross bra"
$values[40,0] = "green l"
$values[41,0] = "cup with"
$values[42,0] = "sport gym"
$values[43,0] = "long bra"
$values[44,0] = "yoga gym"
$values[45,0] = "low back bralette"
$values[46,0] = "longline yoga bra"
$values[47,0] = "criss cross bralette"
$values[48,0] = "sport elastic"
$values[49,0] = "cute back"
$values[50,0] = "pink cross"
$values[51,0] = "green s"
$values[52,0] = "black apparel"
$values[53,0] = "black yoga top"
$values[54,0] = "green m"
$values[55,0] = "yoga bra strappy"
$values[56,0] = "cross back bra"
$values[57,0] = "cute pink"
$values[58,0] = "yoga workout"
$values[59,0] = "sport workout"
$values[60,0] = "cute bralette"
$values[61,0] = "low back top"
$values[62,0] = "criss cross back top"
$values[63,0] = "strappy top"
$values[64,0] = "criss cross front"
$values[65,0] = "strappy sport bra"
$values[66,0] = "cross front bra"
$values[67,0] = "bra elastic"
$values[68,0] = "long line"
$values[69,0] = "sport top"
$values[70,0] = "strappy back bralette"
$values[71,0] = "green sport bra"
$values[72,0] = "bra cup"
$values[73,0] = "wireless sport"
$values[74,0] = "top m"
$values[75,0] = "strappy back top"
$values[76,0] = "gym apparel"
$values[77,0] = "gym back"
$values[78,0] = "back top"
$values[79,0] = "strappy yoga bra"
$values[80,0] = "sport cup"
$values[81,0] = "strappy back bra"
$values[82,0] = "strappy yoga top"
$values[83,0] = "criss cross back"
$values[84,0] = "long l"
$values[85,0] = "low front bra"
$values[86,0] = "gym bra"
$values[87,0] = "workout sport bra"
$values[88,0] = "cross back bralette"
$values[89,0] = "pink criss cross top"
$values[90,0] = "black strappy top"
$values[91,0] = "cross line"
$values[92,0] = "strappy workout bra"
$values[93,0] = "back workout"
$values[94,0] = "yoga bra top"
$values[95,0] = "black strappy"
$values[96,0] = "black criss cross top"
$values[97,0] = "padded yoga bra"
$values[98,0] = "pink bralette top"
$values[99,0] = "green cross"

$ws.Range("A1:A100").Value = $values

# Best-effort: move the active selection to A1 (matches the target sheetView state)
$ws.Range("A1").Select()

